$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "NA" text in C319 (it becomes an empty inline string cell)
$ws.Range("C319").Value = $null

$ws.Range("A320").Value = "2025-11-07"
$ws.Range("B320").Value = "eaux souterraines"
$ws.Range("C320").Value = 64
$ws.Range("D320").Value = 1

$ws.Range("A321").Value = "2025-11-07"
$ws.Range("B321").Value = "ruissellement"
$ws.Range("C321").Value = 66
$ws.Range("D321").Value = 1

$ws.Range("A322").Value = "2025-11-07"
$ws.Range("B322").Value = "eaux souterraines"
$ws.Range("C322").Value = 70
$ws.Range("D322").Value = 2

$ws.Range("A323").Value = "2025-11-07"
$ws.Range("B323").Value = "eaux de surface"
$ws.Range("C323").Value = 70
$ws.Range("D323").Value = 1

$ws.Range("A324").Value = "2025-11-07"
$ws.Range("B324").Value = "eaux souterraines"
$ws.Range("C324").Value = 71
$ws.Range("D324").Value = 1

$ws.Range("A325").Value = "2025-11-07"
$ws.Range("B325").Value = "ruissellement"
$ws.Range("C325").Value = 72
$ws.Range("D325").Value = 1

$ws.Range("A326").Value = "2025-11-07"
$ws.Range("B326").Value = "eaux souterraines"
$ws.Range("C326").Value = 76
$ws.Range("D326").Value = 2

$ws.Range("A327").Value = "2025-11-07"
$ws.Range("B327").Value = "eaux de surface"
$ws.Range("C327").Value = 76
$ws.Range("D327").Value = 1

$ws.Range("A328").Value = "2025-11-07"
$ws.Range("B328").Value = "eaux souterraines"
$ws.Range("C328").Value = 77
$ws.Range("D328").Value = 1

